# chore(login): removed default values
#
# The "Name"/"Email" sign-up sheet shipped with sample/default data
# ("Akshat" / "akshat@gmail.com") pre-filled in the row under the
# headers, including a mailto: hyperlink on the email cell. Remove
# those default values so the form starts out blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the default Name/Email values from row 2.
$dataRow = $ws.Range("A2:B2")
$dataRow.ClearContents()

# The email cell carried a mailto: hyperlink pointing at the removed
# default address - drop it along with the value.
$dataRow.Hyperlinks.Delete()

# Leave the selection on the now-empty entry area.
$ws.Range("A2:C7").Select()
